$d = $word.ActiveDocument

$pairs = @(
    @{old="648×4="; new="256×5="},
    @{old="636×4="; new="849×5="},
    @{old="359×9="; new="703×6="},
    @{old="566×4="; new="917×8="},
    @{old="579×3="; new="423×3="},
    @{old="228×6="; new="243×7="},
    @{old="858×3="; new="401×3="},
    @{old="540×6="; new="156×7="},
    @{old="385×9="; new="972×4="},
    @{old="433×4="; new="136×7="},
    @{old="142×4="; new="673×7="},
    @{old="370×5="; new="536×5="},
    @{old="836×8="; new="497×6="},
    @{old="636×3="; new="874×9="},
    @{old="577×8="; new="674×5="},
    @{old="675×6="; new="669×9="},
    @{old="642×3="; new="461×8="},
    @{old="627×5="; new="567×3="},
    @{old="630×9="; new="787×5="},
    @{old="781×2="; new="131×7="},
    @{old="132×6="; new="887×5="},
    @{old="533×2="; new="756×5="},
    @{old="933×5="; new="276×8="},
    @{old="394×8="; new="956×2="},
    @{old="136×5="; new="520×3="}
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $pair.new, 2)
}
